$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"5.591750551184209E-10"
$ws.Range("E2").Value = [double]"5.591750551184209E-10"

$ws.Range("D3").Value = [double]"8.541227467174458E-05"
$ws.Range("E3").Value = [double]"8.541227467174458E-05"

$ws.Range("D4").Value = [double]"1.524909306636194E-09"
$ws.Range("E4").Value = [double]"1.524909306636194E-09"

$ws.Range("D5").Value = [double]"3.144034209128708E-12"
$ws.Range("E5").Value = [double]"3.144034209128708E-12"

$ws.Range("D6").Value = [double]"1.987192923288383E-08"
$ws.Range("E6").Value = [double]"1.987192923288383E-08"

$ws.Range("D7").Value = 0.9996084498033058
$ws.Range("E7").Value = 0.0003915501966942347

$ws.Range("D8").Value = 0.9999996805554284
$ws.Range("E8").Value = [double]"3.194445715726246E-07"

$ws.Range("D9").Value = 0.9832416654020297
$ws.Range("E9").Value = 0.01675833459797027

$ws.Range("D10").Value = 0.9999999999806397
$ws.Range("E10").Value = [double]"1.93602911480184E-11"

$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.002366774993449494
$ws.Range("E11").Value = 0.9976332250065505
$ws.Range("F11").Value = 0.6063607335090637
$ws.Range("G11").Value = 0.9
